$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Disability" column (column D). This shifts every following
# column one position to the left, updates the dimension, row spans,
# and prunes the now-unused "Disability" shared string automatically.
$ws.Columns("D").Delete()

# Update the scroll/selection state left behind in the worksheet view:
# after deleting column D the user's selection is column D (first of the
# remaining numeric columns) instead of the old column AB.
$null = $ws.Range("D1:D1048576").Select()

# The hidden _xlnm._FilterDatabase defined name still points at the old
# right-most column (AC); shrink it to the new right-most column (AB).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "test_data!_FilterDatabase") {
        $n.RefersTo = "=test_data!`$A`$1:`$AB`$7"
    }
}

# Refresh the persisted custom-sort range (A2:AC7 -> A2:AB7) that is
# stored on the worksheet so it keeps matching the shrunken data range.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$null = $sortObj.SortFields.Add($ws.Range("A1"))
$sortObj.SetRange($ws.Range("A1:AB7"))
$sortObj.Header = 1
$sortObj.Apply()
